$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1936
$ws.Range("K3").Value = 1852
$ws.Range("K4").Value = 397
$ws.Range("K5").Value = 122
$ws.Range("K6").Value = 2371
$ws.Range("K7").Value = 6678

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 188
$ws.Range("K8").Value = 441
$ws.Range("K11").Value = 140
$ws.Range("K19").Value = 186
$ws.Range("K20").Value = 145
$ws.Range("K23").Value = 61
$ws.Range("K27").Value = 75
$ws.Range("K29").Value = 327
$ws.Range("K31").Value = 77
$ws.Range("K33").Value = 269
$ws.Range("K36").Value = 79
$ws.Range("K37").Value = 223
$ws.Range("K42").Value = 235
$ws.Range("K44").Value = 66
$ws.Range("K45").Value = 7
$ws.Range("K47").Value = 41
$ws.Range("K48").Value = 79
$ws.Range("K49").Value = 39
$ws.Range("K53").Value = 102
$ws.Range("K54").Value = 114
$ws.Range("K55").Value = 69
$ws.Range("K63").Value = 24
$ws.Range("K65").Value = 162
$ws.Range("K67").Value = 258
$ws.Range("K69").Value = 18
$ws.Range("K72").Value = 29
$ws.Range("K73").Value = 65
$ws.Range("K76").Value = 97
$ws.Range("K77").Value = 45
$ws.Range("K78").Value = 87
$ws.Range("K83").Value = 146
$ws.Range("K85").Value = 334
$ws.Range("K86").Value = 45
$ws.Range("K88").Value = 88
$ws.Range("K89").Value = 87
$ws.Range("K90").Value = 56
$ws.Range("K92").Value = 30
$ws.Range("K94").Value = 77
$ws.Range("K99").Value = 121
$ws.Range("K101").Value = 6678

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 14
$ws.Range("K6").Value = 44

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 67
$ws.Range("K3").Value = 58
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 45
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 111
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 334

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 128
$ws.Range("K3").Value = 127
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 441

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 56
$ws.Range("K3").Value = 73
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 223

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 14
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 81
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 258

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 21
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 56
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 235

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 6
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 21
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 7
